$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update panel_query_time (column F) timestamps on the "data" sheet
$ws.Range("F2").Value = "2021-10-05 14:33:04.079055"
$ws.Range("F3").Value = "2021-10-05 14:33:04.079064"
$ws.Range("F4").Value = "2021-10-05 14:33:04.079067"
$ws.Range("F5").Value = "2021-10-05 14:33:04.079070"
$ws.Range("F6").Value = "2021-10-05 14:33:04.079073"
$ws.Range("F7").Value = "2021-10-05 14:33:04.079075"
$ws.Range("F8").Value = "2021-10-05 14:33:04.079078"
$ws.Range("F9").Value = "2021-10-05 14:33:04.079081"
$ws.Range("F10").Value = "2021-10-05 14:33:04.079084"
$ws.Range("F11").Value = "2021-10-05 14:33:04.079086"
$ws.Range("F12").Value = "2021-10-05 14:33:04.079089"
$ws.Range("F13").Value = "2021-10-05 14:33:04.079092"
$ws.Range("F14").Value = "2021-10-05 14:33:04.079094"
$ws.Range("F15").Value = "2021-10-05 14:33:04.079097"
$ws.Range("F16").Value = "2021-10-05 14:33:04.079100"
$ws.Range("F17").Value = "2021-10-05 14:33:04.079102"
$ws.Range("F18").Value = "2021-10-05 14:33:04.079105"
$ws.Range("F19").Value = "2021-10-05 14:33:04.079108"
$ws.Range("F20").Value = "2021-10-05 14:33:04.079110"
$ws.Range("F21").Value = "2021-10-05 14:33:04.079113"
$ws.Range("F22").Value = "2021-10-05 14:33:04.079116"
$ws.Range("F23").Value = "2021-10-05 14:33:04.079118"
$ws.Range("F24").Value = "2021-10-05 14:33:04.079121"
$ws.Range("F25").Value = "2021-10-05 14:33:04.079124"
$ws.Range("F26").Value = "2021-10-05 14:33:04.079127"
$ws.Range("F27").Value = "2021-10-05 14:33:04.079130"
$ws.Range("F28").Value = "2021-10-05 14:33:04.079132"
$ws.Range("F29").Value = "2021-10-05 14:33:04.079135"
$ws.Range("F30").Value = "2021-10-05 14:33:04.079138"
$ws.Range("F31").Value = "2021-10-05 14:33:04.079140"
$ws.Range("F32").Value = "2021-10-05 14:33:04.079143"
$ws.Range("F33").Value = "2021-10-05 14:33:04.079145"
$ws.Range("F34").Value = "2021-10-05 14:33:04.079148"
$ws.Range("F35").Value = "2021-10-05 14:33:04.079151"
$ws.Range("F36").Value = "2021-10-05 14:33:04.079154"
$ws.Range("F37").Value = "2021-10-05 14:33:04.079157"
$ws.Range("F38").Value = "2021-10-05 14:33:04.079159"
$ws.Range("F39").Value = "2021-10-05 14:33:04.079162"
$ws.Range("F40").Value = "2021-10-05 14:33:04.079164"
$ws.Range("F41").Value = "2021-10-05 14:33:04.079167"
$ws.Range("F42").Value = "2021-10-05 14:33:04.079170"
$ws.Range("F43").Value = "2021-10-05 14:33:04.079173"
$ws.Range("F44").Value = "2021-10-05 14:33:04.079175"
$ws.Range("F45").Value = "2021-10-05 14:33:04.079178"
$ws.Range("F46").Value = "2021-10-05 14:33:04.079181"
$ws.Range("F47").Value = "2021-10-05 14:33:04.079183"
$ws.Range("F48").Value = "2021-10-05 14:33:04.079186"
$ws.Range("F49").Value = "2021-10-05 14:33:04.079188"
$ws.Range("F50").Value = "2021-10-05 14:33:04.079191"
$ws.Range("F51").Value = "2021-10-05 14:33:04.079194"
$ws.Range("F52").Value = "2021-10-05 14:33:04.079196"
$ws.Range("F53").Value = "2021-10-05 14:33:04.079199"
$ws.Range("F54").Value = "2021-10-05 14:33:04.079202"
$ws.Range("F55").Value = "2021-10-05 14:33:04.079205"
$ws.Range("F56").Value = "2021-10-05 14:33:04.079208"
$ws.Range("F57").Value = "2021-10-05 14:33:04.079215"
$ws.Range("F58").Value = "2021-10-05 14:33:04.079218"
$ws.Range("F59").Value = "2021-10-05 14:33:04.079221"
$ws.Range("F60").Value = "2021-10-05 14:33:04.079223"
$ws.Range("F61").Value = "2021-10-05 14:33:04.079226"
$ws.Range("F62").Value = "2021-10-05 14:33:04.079229"
$ws.Range("F63").Value = "2021-10-05 14:33:04.079231"
$ws.Range("F64").Value = "2021-10-05 14:33:04.079234"
$ws.Range("F65").Value = "2021-10-05 14:33:04.079236"
$ws.Range("F66").Value = "2021-10-05 14:33:04.079240"
$ws.Range("F67").Value = "2021-10-05 14:33:04.079243"
$ws.Range("F68").Value = "2021-10-05 14:33:04.079246"
$ws.Range("F69").Value = "2021-10-05 14:33:04.079248"
$ws.Range("F70").Value = "2021-10-05 14:33:04.079251"
$ws.Range("F71").Value = "2021-10-05 14:33:04.079253"
$ws.Range("F72").Value = "2021-10-05 14:33:04.079256"

# Add a new "metadata" worksheet positioned after "data"
$wsMeta = $wb.Worksheets.Add($null, $ws)
$wsMeta.Name = "metadata"

# Header row
$wsMeta.Range("B1").Value = "data_name"
$wsMeta.Range("C1").Value = "data_id"
$wsMeta.Range("D1").Value = "data_version"
$wsMeta.Range("E1").Value = "data_version_created"
$wsMeta.Range("F1").Value = "panel_query_time"
$wsMeta.Range("G1").Value = "panel_get_request"

# Match header style (bold + border) used by the "data" sheet headers
$ws.Range("B1").Copy()
$wsMeta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$wsMeta.Range("A2").PasteSpecial(-4122)

# Data row
$wsMeta.Range("A2").Value = 0
$wsMeta.Range("B2").Value = "Additional findings_Adult"
$wsMeta.Range("C2").Value = 221
$wsMeta.Range("D2").NumberFormat = "@"
$wsMeta.Range("D2").Value = "0.149"
$wsMeta.Range("D2").Style = "Normal"
$wsMeta.Range("E2").Value = "2021-05-26T08:32:22.955387Z"
$wsMeta.Range("F2").Value = "2021-10-05 14:33:04.075109"
$wsMeta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/221/?format=json"

$wb.Worksheets.Item("data").Activate()
